# Sprint1 daily scrum meeting day3
#
# Updates task status / burn-down tracking on the "Sprint" sheet for day 3
# of the sprint, then leaves the "Chart" sheet as the active/selected tab
# (matching where the user ended up after the meeting).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# "Playing music and flow control" task is now Done; 2 hours were logged
# against it on day 3 (column I).
$ws.Range("F5").Value = "Done"
$ws.Range("I5").Value = 2

# "Find library that implements play functionality" moves to In progress.
$ws.Range("F6").Value = "In progress"

# "Implement basic GUI" moves to In progress; 1 hour logged on day 1
# (column G).
$ws.Range("F9").Value = "In progress"
$ws.Range("G9").Value = 1

# Move the cursor on the Sprint sheet to where the team left off today.
$ws.Range("H8").Select()

# The Chart (burndown) sheet is what's left on screen/active at the end
# of the meeting.
$wsChart = $wb.Worksheets.Item("Chart")
$wsChart.Activate()
